$wb = $excel.ActiveWorkbook

# Give the drone sheets meaningful names instead of the generic DRONE1/2/3
$wb.Worksheets.Item(1).Name = "DJI Phantom 4"
$wb.Worksheets.Item(2).Name = "DJI Mavic 3"
$wb.Worksheets.Item(3).Name = "Custom Drone 1"

# Make the third sheet ("Custom Drone 1") the active/selected tab
$wb.Worksheets.Item(3).Activate()
